# Add a new expense row ("bazar") to the Expenses sheet, and backfill the
# previously-missing "contributions" (I2) cell on the existing row with an
# explicit empty string, matching the target schema (A:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# --- Row 2: add the missing I2 "contributions" cell as an explicit empty string ---
$ws.Range("I2").Formula = '=""'

# --- Row 3: new expense entry ---
$ws.Range("A3").Value = "bazar"
$ws.Range("B3").Value = "Ammu"
$ws.Range("C3").Value = 9.99
# Keep the date as literal text (matches the sheet's existing text-stored dates)
# instead of letting Excel auto-convert it to a date serial number.
$ws.Range("D3").Value = "'2025-12-28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "Hand Cash"
$ws.Range("G3").Value = '["me"]'
$ws.Range("H3").Value = "equal"
$ws.Range("J3").Formula = '=""'
$ws.Range("K3").Value = $false
$ws.Range("L3").Formula = '=""'
$ws.Range("M3").Formula = '=""'
